# The commit only changes the slide order: the slide that used to sit at
# position 11 ("git Grundkonfiguration" / Git Bash basic config) is moved so
# that it now comes right after the slide at position 6, i.e. becomes the
# new slide 7. All other slides keep their relative order (they simply
# shift down by one to make room).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$s.MoveTo(7)
